# Data update from DGS's 2021/09/03 report.
# Appends the newest time-series observation as row 76 of Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 76

# Column A holds a date-like label that is stored as literal text (not a
# real date serial) throughout the sheet, displayed with a "yyyy/mm/dd"
# number format. Temporarily switching the cell to a Text format before
# assigning the value prevents the smart date-parser from converting the
# string into a date serial; restoring the original date number format
# afterwards reproduces the same look (and reuses the existing style) as
# every other cell in the column.
$dateCell = $ws.Cells.Item($row, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2020/09/03"
$dateCell.NumberFormat = "yyyy/mm/dd"

$ws.Cells.Item($row, 2).Value = 295.5
$ws.Cells.Item($row, 3).Value = 302.6
$ws.Cells.Item($row, 4).Value = 0.96
$ws.Cells.Item($row, 5).Value = 0.97
